$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-12-26 Tuesday" "2023-12-27 Wednesday"

Replace-Text "33÷5=" "27÷8="
Replace-Text "57÷3=" "93÷9="
Replace-Text "66÷2=" "30÷3="
Replace-Text "11÷5=" "42÷3="
Replace-Text "23÷7=" "44÷4="

Replace-Text "72÷2=" "72÷5="
Replace-Text "18÷2=" "34÷8="
Replace-Text "35÷4=" "96÷9="
Replace-Text "31÷2=" "83÷3="
Replace-Text "29÷6=" "90÷5="

Replace-Text "25÷6=" "51÷6="
Replace-Text "96÷3=" "92÷5="
Replace-Text "97÷3=" "65÷4="
Replace-Text "38÷2=" "19÷5="
Replace-Text "86÷7=" "87÷6="

Replace-Text "85÷4=" "30÷3="
Replace-Text "27÷4=" "38÷8="
Replace-Text "83÷8=" "38÷7="
Replace-Text "74÷3=" "72÷9="
Replace-Text "12÷9=" "52÷5="

Replace-Text "24÷8=" "25÷8="
Replace-Text "39÷6=" "73÷7="
Replace-Text "59÷5=" "50÷8="
Replace-Text "52÷7=" "37÷7="
Replace-Text "65÷2=" "87÷4="
